$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking
# strings (e.g. "5.22") are stored as text, not auto-converted to numbers,
# matching the inline-string cell type used throughout the sheet.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '62.537.91'
$ws.Range("E2").Value = '  -2.41%  '
$ws.Range("D3").Value = '2.434.75'
$ws.Range("E3").Value = '  -2.64%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = '576.21'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").Value = '143.63'
$ws.Range("E6").Value = '  -4.88%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '2.430.26'
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("E10").Value = '  -5.99%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '5.22'
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("E13").Value = '  -4.15%  '
$ws.Range("D14").Value = '26.37'
$ws.Range("E14").Value = '  -4.14%  '
$ws.Range("E15").Value = '  -5.43%  '
$ws.Range("D16").Value = '2.893.18'
$ws.Range("E16").Value = '  -2.28%  '
$ws.Range("D17").Value = '62.657.22'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("D18").Value = '2.426.34'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("D19").Value = '11.02'
$ws.Range("E19").Value = '  -5.38%  '
$ws.Range("E20").Value = '  -4.53%  '
$ws.Range("D21").Value = '329.57'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("E22").Value = '  -2.46%  '
$ws.Range("D23").Value = '1.99'
$ws.Range("E23").Value = '  -4.35%  '
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").Value = '65.76'
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").Value = '632.38'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '8.91'
$ws.Range("E27").Value = '  +2.10%  '
$ws.Range("D28").Value = '2.565.06'
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = '0.0₃0960'
$ws.Range("E29").Value = '  -9.90%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  -6.80%  '
$ws.Range("D32").Value = '8.04'
$ws.Range("E32").Value = '  -4.77%  '
$ws.Range("E33").Value = '  -2.28%  '
$ws.Range("E34").Value = '  -4.38%  '
$ws.Range("D35").Value = '4.98'
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -5.97%  '
$ws.Range("D38").Value = '0.376'
$ws.Range("E38").Value = '  -2.56%  '
$ws.Range("D39").Value = '18.48'
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").Value = '148.78'
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").Value = '5.28'
$ws.Range("E41").Value = '  -4.94%  '
$ws.Range("E42").Value = '  -5.49%  '
$ws.Range("D43").Value = '42.47'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  -10.31%  '
$ws.Range("E46").Value = '  -4.08%  '
$ws.Range("D47").Value = '3.69'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").Value = '0.0523'
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("D49").Value = '0.596'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("D50").Value = '19.64'
$ws.Range("E50").Value = '  -7.61%  '
$ws.Range("D51").Value = '0.0₆0233'
$ws.Range("E51").Value = '  +4.60%  '

# Restore default (General) formatting now that the text values are in place,
# so the cells do not carry a leftover style index.
$dataRange.ClearFormats()

